$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1341219.2
$ws.Range("I17").Value = 1030
$ws.Range("K17").Value = 3090
$ws.Range("M17").Value = -2922
$ws.Range("H19").Value = 1695.375
$ws.Range("I19").Value = 1665
$ws.Range("K19").Value = 1665
$ws.Range("M19").Value = -1490
$ws.Range("H58").Value = 1969.2222
$ws.Range("I58").Value = 2747.5
$ws.Range("J58").Value = 1346.6
$ws.Range("K58").Value = 8242.5
$ws.Range("L58").Value = 4039.8
$ws.Range("M58").Value = -8092.5
$ws.Range("N58").Value = -4339.799999999999
$ws.Range("H70").Value = 41671496
$ws.Range("I70").Value = 994.5
$ws.Range("K70").Value = 2983.5
$ws.Range("M70").Value = -2713.5
$ws.Range("H73").Value = 41671496
$ws.Range("I73").Value = 994.5
$ws.Range("K73").Value = 2983.5
$ws.Range("M73").Value = -2047.5
$ws.Range("H106").Value = 3930964.5
$ws.Range("I106").Value = 4771174
$ws.Range("J106").Value = 9988.333000000001
$ws.Range("K106").Value = 4771174
$ws.Range("L106").Value = 9988.333000000001
$ws.Range("M106").Value = -4770543
$ws.Range("N106").Value = -11250.333
$ws.Range("H111").Value = 1090
$ws.Range("I111").Value = 1022.8889
$ws.Range("J111").Value = 1331.6
$ws.Range("K111").Value = 3068.6667
$ws.Range("L111").Value = 3994.8
$ws.Range("M111").Value = -1.666700000000219
$ws.Range("N111").Value = -10128.8
$ws.Range("H116").Value = 29939.375
$ws.Range("I116").Value = 33073.57
$ws.Range("J116").Value = 8000
$ws.Range("K116").Value = 33073.57
$ws.Range("L116").Value = 8000
$ws.Range("M116").Value = -29631.57
$ws.Range("N116").Value = -14884
$ws.Range("H132").Value = 4720.758
$ws.Range("I132").Value = 4926.1665
$ws.Range("J132").Value = 2666.6667
$ws.Range("K132").Value = 14778.4995
$ws.Range("L132").Value = 8000.000100000001
$ws.Range("M132").Value = -12248.4995
$ws.Range("N132").Value = -13060.0001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2095.8948
$ws.Range("I2").Value = 1110.5834
$ws.Range("J2").Value = 3785
$ws.Range("K2").Value = 1110.5834
$ws.Range("L2").Value = 3785
$ws.Range("M2").Value = -997.5834
$ws.Range("N2").Value = -4011
$ws.Range("H32").Value = 4484.0625
$ws.Range("I32").Value = 4099.2544
$ws.Range("K32").Value = 4099.2544
$ws.Range("M32").Value = -3812.2544
$ws.Range("H45").Value = 3917.35
$ws.Range("I45").Value = 3426.5334
$ws.Range("K45").Value = 3426.5334
$ws.Range("M45").Value = -3049.5334
$ws.Range("H97").Value = 1626.5
$ws.Range("I97").Value = 1418.3334
$ws.Range("K97").Value = 1418.3334
$ws.Range("M97").Value = -922.3334
$ws.Range("H102").Value = 2900
$ws.Range("I102").Value = 2300
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 2300
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = -678
$ws.Range("N102").Value = -6744
$ws.Range("H106").Value = 104974.5
$ws.Range("J106").Value = 104974.5
$ws.Range("L106").Value = 104974.5
$ws.Range("N106").Value = -107498.5
$ws.Range("H110").Value = 1601.8125
$ws.Range("I110").Value = 1475.2667
$ws.Range("K110").Value = 1475.2667
$ws.Range("M110").Value = 569.7333000000001
$ws.Range("H116").Value = 2095.8948
$ws.Range("I116").Value = 1110.5834
$ws.Range("J116").Value = 3785
$ws.Range("K116").Value = 1110.5834
$ws.Range("L116").Value = 3785
$ws.Range("M116").Value = 1183.4166
$ws.Range("N116").Value = -8373
$ws.Range("H122").Value = 4412.343
$ws.Range("I122").Value = 2809.7917
$ws.Range("J122").Value = 7908.8184
$ws.Range("K122").Value = 8429.375100000001
$ws.Range("L122").Value = 23726.4552
$ws.Range("M122").Value = -5979.375100000001
$ws.Range("N122").Value = -28626.4552

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2095.8948
$ws.Range("I3").Value = 1110.5834
$ws.Range("J3").Value = 3785
$ws.Range("K3").Value = 1110.5834
$ws.Range("L3").Value = 3785
$ws.Range("M3").Value = -996.5834
$ws.Range("N3").Value = -4013
$ws.Range("H10").Value = 2751
$ws.Range("I10").Value = 2001.3334
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 2001.3334
$ws.Range("L10").Value = 5000
$ws.Range("M10").Value = -1861.3334
$ws.Range("N10").Value = -5280
$ws.Range("H86").Value = 869.1818
$ws.Range("I86").Value = 745.5
$ws.Range("J86").Value = 1199
$ws.Range("K86").Value = 745.5
$ws.Range("L86").Value = 1199
$ws.Range("M86").Value = 377.5
$ws.Range("N86").Value = -3445
$ws.Range("H89").Value = 869.1818
$ws.Range("I89").Value = 745.5
$ws.Range("J89").Value = 1199
$ws.Range("K89").Value = 3727.5
$ws.Range("L89").Value = 5995
$ws.Range("M89").Value = 1888.5
$ws.Range("N89").Value = -17227
$ws.Range("H94").Value = 2831.3333
$ws.Range("I94").Value = 2641.1428
$ws.Range("J94").Value = 3497
$ws.Range("K94").Value = 2641.1428
$ws.Range("L94").Value = 3497
$ws.Range("M94").Value = -2190.1428
$ws.Range("N94").Value = -4399
$ws.Range("H134").Value = 2187.818
$ws.Range("I134").Value = 2101.2
$ws.Range("K134").Value = 6303.599999999999
$ws.Range("M134").Value = -3768.599999999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 74
$ws.Range("J7").Value = 65.333336
$ws.Range("L7").Value = 65.333336
$ws.Range("N7").Value = -291.333336
$ws.Range("H31").Value = 3305.204
$ws.Range("I31").Value = 2159.65
$ws.Range("K31").Value = 2159.65
$ws.Range("M31").Value = -1864.65
$ws.Range("H34").Value = 3305.204
$ws.Range("I34").Value = 2159.65
$ws.Range("K34").Value = 2159.65
$ws.Range("M34").Value = -1957.65
$ws.Range("H105").Value = 1971
$ws.Range("I105").Value = 1793.0555
$ws.Range("J105").Value = 2504.8333
$ws.Range("K105").Value = 1793.0555
$ws.Range("L105").Value = 2504.8333
$ws.Range("M105").Value = -46.05549999999994
$ws.Range("N105").Value = -5998.8333
$ws.Range("H110").Value = 66322.664
$ws.Range("J110").Value = 66322.664
$ws.Range("L110").Value = 66322.664
$ws.Range("N110").Value = -74502.664

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 104174030
$ws.Range("I4").Value = 93050720
$ws.Range("J4").Value = 122998080
$ws.Range("K4").Value = 279152160
$ws.Range("L4").Value = 368994240
$ws.Range("M4").Value = -279152048
$ws.Range("N4").Value = -368994464
$ws.Range("H140").Value = 1550.7188
$ws.Range("I140").Value = 1260.3334
$ws.Range("J140").Value = 3118.8
$ws.Range("K140").Value = 3781.0002
$ws.Range("L140").Value = 9356.400000000001
$ws.Range("M140").Value = 1398.9998
$ws.Range("N140").Value = -19716.4

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3710.25
$ws.Range("I132").Value = 3413.6667
$ws.Range("K132").Value = 10241.0001
$ws.Range("M132").Value = -7711.000100000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2121.7778
$ws.Range("I16").Value = 2121.7778
$ws.Range("K16").Value = 2121.7778
$ws.Range("M16").Value = -1951.7778
$ws.Range("H18").Value = 851334.2
$ws.Range("I18").Value = 29005
$ws.Range("K18").Value = 29005
$ws.Range("M18").Value = -28833
$ws.Range("H61").Value = 2123.1428
$ws.Range("I61").Value = 2735.9167
$ws.Range("J61").Value = 1306.1111
$ws.Range("K61").Value = 2735.9167
$ws.Range("L61").Value = 1306.1111
$ws.Range("M61").Value = -2533.9167
$ws.Range("N61").Value = -1710.1111
$ws.Range("H82").Value = 17376.846
$ws.Range("I82").Value = 2999.3333
$ws.Range("J82").Value = 21690.1
$ws.Range("K82").Value = 2999.3333
$ws.Range("L82").Value = 21690.1
$ws.Range("M82").Value = -2638.3333
$ws.Range("N82").Value = -22412.1
$ws.Range("H85").Value = 17376.846
$ws.Range("I85").Value = 2999.3333
$ws.Range("J85").Value = 21690.1
$ws.Range("K85").Value = 2999.3333
$ws.Range("L85").Value = 21690.1
$ws.Range("M85").Value = -1751.3333
$ws.Range("N85").Value = -24186.1
$ws.Range("H99").Value = 81293.164
$ws.Range("J99").Value = 99439.75
$ws.Range("L99").Value = 99439.75
$ws.Range("N99").Value = -105429.75
$ws.Range("H102").Value = 11985
$ws.Range("J102").Value = 11985
$ws.Range("L102").Value = 11985
$ws.Range("N102").Value = -18475
$ws.Range("H113").Value = 2123.1428
$ws.Range("I113").Value = 2735.9167
$ws.Range("J113").Value = 1306.1111
$ws.Range("K113").Value = 2735.9167
$ws.Range("L113").Value = 1306.1111
$ws.Range("M113").Value = -565.9167000000002
$ws.Range("N113").Value = -5646.1111
$ws.Range("H132").Value = 2295.8572
$ws.Range("I132").Value = 1664.7
$ws.Range("J132").Value = 3873.75
$ws.Range("K132").Value = 4994.1
$ws.Range("L132").Value = 11621.25
$ws.Range("M132").Value = -2464.1
$ws.Range("N132").Value = -16681.25

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5100.3
$ws.Range("I62").Value = 3875
$ws.Range("J62").Value = 5917.1665
$ws.Range("K62").Value = 3875
$ws.Range("L62").Value = 5917.1665
$ws.Range("M62").Value = -3251
$ws.Range("N62").Value = -7165.1665
$ws.Range("H65").Value = 5100.3
$ws.Range("I65").Value = 3875
$ws.Range("J65").Value = 5917.1665
$ws.Range("K65").Value = 19375
$ws.Range("L65").Value = 29585.8325
$ws.Range("M65").Value = -16255
$ws.Range("N65").Value = -35825.8325
$ws.Range("H96").Value = 11064.3
$ws.Range("I96").Value = 7891.625
$ws.Range("J96").Value = 13179.417
$ws.Range("K96").Value = 7891.625
$ws.Range("L96").Value = 13179.417
$ws.Range("M96").Value = -6518.625
$ws.Range("N96").Value = -15925.417
$ws.Range("H107").Value = 487.7647
$ws.Range("J107").Value = 472.2
$ws.Range("L107").Value = 1416.6
$ws.Range("N107").Value = -5256.6
$ws.Range("H122").Value = 4697.2593
$ws.Range("I122").Value = 2683.65
$ws.Range("J122").Value = 10450.429
$ws.Range("K122").Value = 8050.950000000001
$ws.Range("L122").Value = 31351.287
$ws.Range("M122").Value = -5600.950000000001
$ws.Range("N122").Value = -36251.287
